$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new header cells P1, Q1 (copy format from O1 so they get style index 1) ---
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Add new columns P, Q (value 0) for rows 2-25 ---
$ws.Range("P2:Q25").Value = 0

# --- Update changed numeric values in columns B, C, D, F, G, I for rows 2-25 ---
$ws.Range("B2").Value = 3.449428331641002
$ws.Range("C2").Value = 0.9982227907523225
$ws.Range("D2").Value = 0.03172138027976246
$ws.Range("F2").Value = 1.814263910620483
$ws.Range("G2").Value = 0.0007938166815006416
$ws.Range("I2").Value = 0.002040037396629835

$ws.Range("B3").Value = 2.993295657310057
$ws.Range("C3").Value = 0.8637839152792992
$ws.Range("D3").Value = 0.03239356478296429
$ws.Range("F3").Value = 1.624637643376602
$ws.Range("G3").Value = 0.0008009136028366262
$ws.Range("I3").Value = 0.0005415994938058333

$ws.Range("B4").Value = 2.714879073436464
$ws.Range("C4").Value = 0.7827315438028677
$ws.Range("D4").Value = 0.03276040081950882
$ws.Range("F4").Value = 1.509036365549221
$ws.Range("G4").Value = 0.0008053941136406356
$ws.Range("I4").Value = 0.0004733647622363613

$ws.Range("B5").Value = 2.601712334019965
$ws.Range("C5").Value = 0.7515755768533836
$ws.Range("D5").Value = 0.03271940708173204
$ws.Range("F5").Value = 1.455012734718579
$ws.Range("G5").Value = 0.0008072665837825106
$ws.Range("I5").Value = 0.0007113638850162474

$ws.Range("B6").Value = 2.582883763513337
$ws.Range("C6").Value = 0.7483124265440608
$ws.Range("D6").Value = 0.03249804333232653
$ws.Range("F6").Value = 1.437445664522158
$ws.Range("G6").Value = 0.000807597168562532
$ws.Range("I6").Value = 0.0008540287709379868

$ws.Range("B7").Value = 2.713187553185435
$ws.Range("C7").Value = 0.7874845422498709
$ws.Range("D7").Value = 0.03216823667292878
$ws.Range("F7").Value = 1.484705568971577
$ws.Range("G7").Value = 0.0008054670675243055
$ws.Range("I7").Value = 0.0007047959113215541

$ws.Range("B8").Value = 3.291485219529761
$ws.Range("C8").Value = 0.9584315267965735
$ws.Range("D8").Value = 0.03116883041336749
$ws.Range("F8").Value = 1.717098490305418
$ws.Range("G8").Value = 0.0007962998941470674
$ws.Range("I8").Value = 0.001581273689130747

$ws.Range("B9").Value = 4.443011327753197
$ws.Range("C9").Value = 1.298344364862317
$ws.Range("D9").Value = 0.02993710216353129
$ws.Range("F9").Value = 2.225773536441793
$ws.Range("G9").Value = 0.0007791603013742216
$ws.Range("I9").Value = 0.01062721742675876

$ws.Range("B10").Value = 5.305526218380351
$ws.Range("C10").Value = 1.56266870637171
$ws.Range("D10").Value = 0.02854187158562738
$ws.Range("F10").Value = 2.599078181730164
$ws.Range("G10").Value = 0.0007670609992566027
$ws.Range("I10").Value = 0.02332236619447237

$ws.Range("B11").Value = 5.702463316671128
$ws.Range("C11").Value = 1.695150015529009
$ws.Range("D11").Value = 0.02690158284183397
$ws.Range("F11").Value = 2.732620438162456
$ws.Range("G11").Value = 0.0007617051299523281
$ws.Range("I11").Value = 0.03092685545455609

$ws.Range("B12").Value = 5.853976995486789
$ws.Range("C12").Value = 1.74125032028121
$ws.Range("D12").Value = 0.02678925719047065
$ws.Range("F12").Value = 2.805121269063136
$ws.Range("G12").Value = 0.0007596497882525117
$ws.Range("I12").Value = 0.03395731172641359

$ws.Range("B13").Value = 5.821369806137852
$ws.Range("C13").Value = 1.730355516662257
$ws.Range("D13").Value = 0.02691548642701136
$ws.Range("F13").Value = 2.7937567358195
$ws.Range("G13").Value = 0.0007600847718871556
$ws.Range("I13").Value = 0.03326810816659798

$ws.Range("B14").Value = 5.714937202006524
$ws.Range("C14").Value = 1.698531770833597
$ws.Range("D14").Value = 0.02693544161671468
$ws.Range("F14").Value = 2.740386170866088
$ws.Range("G14").Value = 0.0007615327611190859
$ws.Range("I14").Value = 0.03115985350385309

$ws.Range("B15").Value = 5.649740136347646
$ws.Range("C15").Value = 1.680965767896112
$ws.Range("D15").Value = 0.02674985480474312
$ws.Range("F15").Value = 2.699396207634749
$ws.Range("G15").Value = 0.0007624352174286817
$ws.Range("I15").Value = 0.0299580457551416

$ws.Range("B16").Value = 5.278697796008942
$ws.Range("C16").Value = 1.569605084833768
$ws.Range("D16").Value = 0.02695271143622691
$ws.Range("F16").Value = 2.519720400877446
$ws.Range("G16").Value = 0.0007675348952299834
$ws.Range("I16").Value = 0.0232828947889967

$ws.Range("B17").Value = 5.052524724942089
$ws.Range("C17").Value = 1.50215881308543
$ws.Range("D17").Value = 0.02708162223982136
$ws.Range("F17").Value = 2.41098210759165
$ws.Range("G17").Value = 0.0007706821269455662
$ws.Range("I17").Value = 0.01963970417161942

$ws.Range("B18").Value = 4.923200573215865
$ws.Range("C18").Value = 1.458954490730605
$ws.Range("D18").Value = 0.02764265834598412
$ws.Range("F18").Value = 2.370087421251597
$ws.Range("G18").Value = 0.0007724629092008239
$ws.Range("I18").Value = 0.01752388413944761

$ws.Range("B19").Value = 4.879276451817987
$ws.Range("C19").Value = 1.448003433303086
$ws.Range("D19").Value = 0.02745578895909873
$ws.Range("F19").Value = 2.339989425491353
$ws.Range("G19").Value = 0.0007730957371923109
$ws.Range("I19").Value = 0.01696520041713789

$ws.Range("B20").Value = 5.076572415876285
$ws.Range("C20").Value = 1.508937579475855
$ws.Range("D20").Value = 0.02710573035128583
$ws.Range("F20").Value = 2.424166465931975
$ws.Range("G20").Value = 0.0007703434467931631
$ws.Range("I20").Value = 0.01999738226909642

$ws.Range("B21").Value = 5.745806434683118
$ws.Range("C21").Value = 1.713135852033304
$ws.Range("D21").Value = 0.02636761234774454
$ws.Range("F21").Value = 2.732090373252191
$ws.Range("G21").Value = 0.0007611479446893474
$ws.Range("I21").Value = 0.03192115324456424

$ws.Range("B22").Value = 6.189017870973885
$ws.Range("C22").Value = 1.843065342267494
$ws.Range("D22").Value = 0.02661427810312667
$ws.Range("F22").Value = 2.969222859553071
$ws.Range("G22").Value = 0.0007551392189640411
$ws.Range("I22").Value = 0.04113182192123954

$ws.Range("B23").Value = 5.952378374591888
$ws.Range("C23").Value = 1.767049654751872
$ws.Range("D23").Value = 0.02717120523743866
$ws.Range("F23").Value = 2.871098226769163
$ws.Range("G23").Value = 0.000758292914541373
$ws.Range("I23").Value = 0.03589869942773749

$ws.Range("B24").Value = 5.066309906723802
$ws.Range("C24").Value = 1.495922873222128
$ws.Range("D24").Value = 0.02815156441209155
$ws.Range("F24").Value = 2.463219505788445
$ws.Range("G24").Value = 0.0007704181125877175
$ws.Range("I24").Value = 0.01951105419918786

$ws.Range("B25").Value = 4.128569216504843
$ws.Range("C25").Value = 1.214108523772381
$ws.Range("D25").Value = 0.02924206279023522
$ws.Range("F25").Value = 2.043038658649309
$ws.Range("G25").Value = 0.000783786542006243
$ws.Range("I25").Value = 0.007537696458468801
